$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '36.331.52'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -0.04%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.927.14'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -2.70%  '
$ws.Range('E4').Value = '  +0.01%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '240.07'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -2.02%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.604'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -3.40%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -0.12%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '55.80'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  -5.52%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.355'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -4.88%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.0825'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('E11').Value = '  -2.26%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '2.208.66'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -2.76%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '20.67'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -10.85%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '0.787'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -8.38%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '13.15'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -5.69%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '5.07'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -6.84%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '1.933.07'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -2.55%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '36.212.95'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +0.12%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '68.29'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -2.86%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0851'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -2.84%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '225.36'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -3.58%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '4.87'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -7.88%  '
$ws.Range('E23').Value = '  +0.10%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '2.30'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -8.84%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.25'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -1.83%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '9.01'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -8.97%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '159.66'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('E28').Value = '  -3.57%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '18.99'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -3.92%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.116'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -3.20%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.08'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -7.63%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '4.46'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -8.28%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.0611'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -10.83%  '
$ws.Range('E34').Value = '  +0.01%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '4.08'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -7.14%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '5.94'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -4.21%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.78'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  -1.69%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '2.11'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -5.84%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '2.93'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -1.26%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0954'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  -1.85%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.0206'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('E43').Value = '  -8.54%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '15.33'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -4.85%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.317.48'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -3.36%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -7.88%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '83.91'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -8.75%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '6.92'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -6.99%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '2.81'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -0.41%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '2.099.74'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -2.75%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '42.55'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -5.34%  '
